$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6557395458221436
$ws.Range("B1").Value = 0.521864652633667
$ws.Range("C1").Value = 4.860641956329346
$ws.Range("D1").Value = 2.812721729278564
$ws.Range("E1").Value = 1.215299010276794
